$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base roster of 8 workers (doc number, name), in the original order.
# The account statement previously listed each worker's two mora periods
# (1608 then 1607) back to back; the update regroups the data so all
# "1607" period rows come first (in worker order) followed by all "1608"
# period rows (same worker order) -- i.e. sorted/grouped by Periodo Mora.
$workers = @(
    @{ Doc = "45554102";   Name = "VICELYS JULIO JULIO" },
    @{ Doc = "45483945";   Name = "DORIS MAGALY OROZCO BOSSIO" },
    @{ Doc = "73192259";   Name = "VICTOR MANUEL PADILLA PAREJA" },
    @{ Doc = "1002189531"; Name = "INDRINA NAYARITH FERRER DE AVILA" },
    @{ Doc = "1047492034"; Name = "CARLOS YESID LEMUS MONTOYA" },
    @{ Doc = "73138664";   Name = "LIZARDO ENRIQUE VILLALBA BARBOZA" },
    @{ Doc = "73115645";   Name = "ROBERT ALFONSO COCHERO LAMBIS" },
    @{ Doc = "1047382051"; Name = "FELICIDAD EXTREMOR PADILLA" }
)

$periods = @("1607", "1608")

$row = 16
foreach ($period in $periods) {
    foreach ($w in $workers) {
        $ws.Cells.Item($row, 2).Value = "CC"
        $ws.Cells.Item($row, 3).Value = $w.Doc
        $ws.Cells.Item($row, 4).Value = $w.Name
        $ws.Cells.Item($row, 5).Value = $period
        $ws.Cells.Item($row, 6).Value = 25774
        # Corrected Salario Basico -- the previous 1500000 entries for
        # LIZARDO ENRIQUE VILLALBA BARBOZA were wrong; every worker earns
        # the same base salary of 644350.
        $ws.Cells.Item($row, 7).Value = 644350
        $row = $row + 1
    }
}
